$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "symbol" column (A2:A9) and the value grid (B2:H9) ---
# New symbol labels (shared strings change from WL/BH/BM/... to SC/A/B/C/D/E/F/G)
$symbols = @("SC", "A", "B", "C", "D", "E", "F", "G")
for ($i = 0; $i -lt $symbols.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $symbols[$i]
}

# New numeric grid for rows 2-9, columns B..H (cols 2..8)
$grid = @(
    @(0, 0, 0, 1, 0, 0, 0),
    @(24, 0, 24, 0, 0, 24, 0),
    @(0, 24, 0, 0, 24, 0, 24),
    @(0, 24, 0, 24, 0, 24, 0),
    @(24, 0, 0, 24, 0, 24, 0),
    @(0, 0, 24, 0, 24, 0, 24),
    @(0, 24, 0, 24, 0, 0, 24),
    @(24, 0, 24, 0, 24, 0, 0)
)
for ($i = 0; $i -lt $grid.Length; $i++) {
    $r = 2 + $i
    $vals = $grid[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $c = 2 + $j
        $ws.Cells.Item($r, $c).Value = $vals[$j]
    }
}

# --- Remove the old rows 11-15 (their data is no longer part of the table) ---
$ws.Rows("11:15").Delete()

# --- Row 10 becomes an empty "next row" placeholder with a red outline border ---
$row10 = $ws.Range("A10:H10")
$row10.ClearContents()
$row10.ClearFormats()
$row10.Interior.ColorIndex = 9
$row10.BorderAround(1, 2, 10)
